$wb = $excel.ActiveWorkbook

function Set-Row {
    param($ws, [int]$row, [object[]]$values)
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

# ---------------------------------------------------------------------------
# Sheet "Home win" -> A1:F4
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Home win")

Set-Row $ws1 2 @("18-12-2024 14:30","PORTUGAL","LIGA REVELAÇÃO U23","Torreense U23 - Sporting Braga U23",80,2.1)
Set-Row $ws1 3 @("18-12-2024 10:00","TURKEY","CUP","Keçiörengücü - Sivasspor",70,3)
Set-Row $ws1 4 @("19-12-2024 18:00","SPAIN","SEGUNDA DIVISIÓN","Cadiz - Burgos",70,1.91)

# ---------------------------------------------------------------------------
# Sheet "Draw" -> A1:F4 (new column F "Draw Odds")
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Draw")

$ws2.Range("E1").Copy()
$ws2.Cells.Item(1, 6).PasteSpecial(-4122)  # xlPasteFormats
$ws2.Cells.Item(1, 6).Value = "Draw Odds"

Set-Row $ws2 2 @("18-12-2024 14:00","CAMEROON","ELITE ONE","Stade Renard - Canon",65,2.8)
$ws2.Cells.Item(3, 6).Value = 3.6
Set-Row $ws2 4 @("19-12-2024 18:00","ROMANIA","CUPA ROMÂNIEI","CS Afumati - Arges Pitesti",66.7,3.1)

# ---------------------------------------------------------------------------
# Sheet "Btts" -> A1:F10 (was A1:F11 - row 11 removed, content reshuffled)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Btts")

Set-Row $ws3 2  @("18-12-2024 17:30","ITALY","COPPA ITALIA","Atalanta - Cesena",80,2.1)
Set-Row $ws3 3  @("18-12-2024 19:00","NETHERLANDS","KNVB BEKER","ASWH - Heerenveen",76,1.7)
Set-Row $ws3 4  @("18-12-2024 19:00","NETHERLANDS","KNVB BEKER","AFC Amsterdam - Utrecht",76,1.91)
Set-Row $ws3 5  @("18-12-2024 15:00","ROMANIA","CUPA ROMÂNIEI","Politehnica Iasi - AFC Hermannstadt",80,1.91)
Set-Row $ws3 6  @("18-12-2024 18:00","SPAIN","SEGUNDA DIVISIÓN","Racing Ferrol - Almeria",76,1.77)
Set-Row $ws3 7  @("24-11-2024 16:00","SPAIN","SEGUNDA DIVISIÓN RFEF - GROUP 1","Bergantiños - Deportivo La Coruña II",80,1.8)
Set-Row $ws3 8  @("19-12-2024 20:00","WORLD","UEFA EUROPA CONFERENCE LEAGUE","Larne - Gent",76,2.1)
Set-Row $ws3 9  @("19-12-2024 20:00","WORLD","UEFA EUROPA CONFERENCE LEAGUE","Celje - The New Saints",76,1.8)
Set-Row $ws3 10 @("19-12-2024 20:15","SPAIN","SEGUNDA DIVISIÓN","Huesca - Tenerife",76.7,2.37)

$ws3.Rows.Item(11).Delete()

# ---------------------------------------------------------------------------
# Sheet "Over_Under" -> A1:H5 (was A1:H6 - row 6 removed, content reshuffled)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Over_Under")

Set-Row $ws4 2 @("18-12-2024 10:00","TURKEY","CUP","Keçiörengücü - Sivasspor",80,1.8,60,3)
Set-Row $ws4 3 @("19-12-2024 20:00","WORLD","UEFA EUROPA CONFERENCE LEAGUE","1. FC Heidenheim - FC ST. Gallen",70,1.53,60,2.25)
Set-Row $ws4 4 @("19-12-2024 20:00","WORLD","UEFA EUROPA CONFERENCE LEAGUE","Djurgardens IF - Legia Warszawa",100,1.85,40,3.25)
Set-Row $ws4 5 @("19-12-2024 20:00","WORLD","UEFA EUROPA CONFERENCE LEAGUE","TSC Backa Topola - FC Noah",86.7,1.65,60,2.6)

$ws4.Rows.Item(6).Delete()
